# Bug fix cost calculation / labour cost
# Corrects the operational_rate_energy_[€/kWh] value (column G, row 3) on the
# shredding, extrusion and conditioning sheets from 0.0439 to 0.207.
# Also restores the view/selection state left behind by the author while
# making the fix (active cell, zoom, tab selection).

$wb = $excel.ActiveWorkbook

$shredding    = $wb.Worksheets.Item(1)   # "shredding"
$extrusion    = $wb.Worksheets.Item(2)   # "extrusion"
$granulate    = $wb.Worksheets.Item(3)   # "granulate"
$conditioning = $wb.Worksheets.Item(4)   # "conditioning"

# --- Bug fix: correct the energy operational rate value used in the labour /
#     machine cost calculation on the affected sheets ---
$shredding.Range("G3").Value = 0.207
$extrusion.Range("G3").Value = 0.207
$conditioning.Range("G3").Value = 0.207

# --- Restore view / selection state for each sheet ---

# shredding: zoom stays the same, active cell moves to G4
[void]$shredding.Range("G4").Select()

# extrusion: scroll so column C is the left-most visible column, active cell G3
$extrusion.Application.ActiveWindow.ScrollColumn = 3
[void]$extrusion.Range("G3").Select()

# conditioning: zoom changes from 86 to 125, active cell moves to H7
[void]$conditioning.Range("H7").Select()
$conditioning.Application.ActiveWindow.Zoom = 125

# granulate becomes the active sheet/tab, selection unchanged (F6)
[void]$granulate.Activate()
